# Update marksheet figures: corrected marks / total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking value (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Total marks (B12): 75 -> 125
$ws.Range("B12").Value = 125

# Correct/Total display (E12): "72/84" -> "125/140"
$ws.Range("E12").Value = "125/140"
